$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Cells.Item(33, 8).Value = 873.6  # H33: 304.8889 -> 873.6
$ws.Cells.Item(33, 9).Value = 1576.4  # I33: 470 -> 1576.4
$ws.Cells.Item(33, 10).Value = 170.8  # J33: 172.8 -> 170.8
$ws.Cells.Item(33, 11).Value = 1576.4  # K33: 470 -> 1576.4
$ws.Cells.Item(33, 12).Value = 170.8  # L33: 172.8 -> 170.8
$ws.Cells.Item(33, 13).Value = -1347.4  # M33: -241 -> -1347.4
$ws.Cells.Item(33, 14).Value = -628.8  # N33: -630.8 -> -628.8
# Row 112
$ws.Cells.Item(112, 8).Value = 2746.4792  # H112: 4226.869 -> 2746.4792
$ws.Cells.Item(112, 9).Value = 940  # I112: 980 -> 940
$ws.Cells.Item(112, 10).Value = 2784.9148  # J112: 4306.061 -> 2784.9148
$ws.Cells.Item(112, 11).Value = 2820  # K112: 2940 -> 2820
$ws.Cells.Item(112, 12).Value = 8354.7444  # L112: 12918.183 -> 8354.7444
$ws.Cells.Item(112, 13).Value = -1712  # M112: -1832 -> -1712
$ws.Cells.Item(112, 14).Value = -10570.7444  # N112: -15134.183 -> -10570.7444
# Row 125
$ws.Cells.Item(125, 8).Value = 17858940  # H125: 11364814 -> 17858940
$ws.Cells.Item(125, 9).Value = 62500700  # I125: 31250694 -> 62500700
$ws.Cells.Item(125, 10).Value = 2236  # J125: 1454.2858 -> 2236
$ws.Cells.Item(125, 11).Value = 562506300  # K125: 281256246 -> 562506300
$ws.Cells.Item(125, 12).Value = 20124  # L125: 13088.5722 -> 20124
$ws.Cells.Item(125, 13).Value = -562503840  # M125: -281253786 -> -562503840
$ws.Cells.Item(125, 14).Value = -25044  # N125: -18008.5722 -> -25044

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Cells.Item(122, 8).Value = 501000  # H122: 168425 -> 501000
$ws.Cells.Item(122, 9).Value = 501000  # I122: 251250 -> 501000
$ws.Cells.Item(122, 10).Value = 0  # J122: 2775 -> 0
$ws.Cells.Item(122, 11).Value = 1503000  # K122: 753750 -> 1503000
$ws.Cells.Item(122, 12).Value = 0  # L122: 8325 -> 0
$ws.Cells.Item(122, 13).Value = -1500550  # M122: -751300 -> -1500550
$ws.Cells.Item(122, 14).ClearContents()  # N122: -13225 -> (removed)
# Row 140
$ws.Cells.Item(140, 8).Value = 96803.5  # H140: 78605.17999999999 -> 96803.5
$ws.Cells.Item(140, 10).Value = 96803.5  # J140: 78605.17999999999 -> 96803.5
$ws.Cells.Item(140, 12).Value = 96803.5  # L140: 78605.17999999999 -> 96803.5
$ws.Cells.Item(140, 14).Value = -107163.5  # N140: -88965.17999999999 -> -107163.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Cells.Item(25, 8).Value = 64277.332  # H25: 40944 -> 64277.332
$ws.Cells.Item(25, 10).Value = 95016  # J25: 60016 -> 95016
$ws.Cells.Item(25, 12).Value = 95016  # L25: 60016 -> 95016
$ws.Cells.Item(25, 14).Value = -95486  # N25: -60486 -> -95486
# Row 107
$ws.Cells.Item(107, 8).Value = 168851.83  # H107: 54402.473 -> 168851.83
$ws.Cells.Item(107, 9).Value = 251777.75  # I107: 64146.688 -> 251777.75
$ws.Cells.Item(107, 10).Value = 3000  # J107: 2433.3333 -> 3000
$ws.Cells.Item(107, 11).Value = 251777.75  # K107: 64146.688 -> 251777.75
$ws.Cells.Item(107, 12).Value = 3000  # L107: 2433.3333 -> 3000
$ws.Cells.Item(107, 13).Value = -249857.75  # M107: -62226.688 -> -249857.75
$ws.Cells.Item(107, 14).Value = -6840  # N107: -6273.3333 -> -6840
# Row 140
$ws.Cells.Item(140, 8).Value = 59650  # H140: 51256 -> 59650
$ws.Cells.Item(140, 10).Value = 59650  # J140: 51256 -> 59650
$ws.Cells.Item(140, 12).Value = 59650  # L140: 51256 -> 59650
$ws.Cells.Item(140, 14).Value = -70010  # N140: -61616 -> -70010

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5618.5557  # H31: 7626.9814 -> 5618.5557
$ws.Cells.Item(31, 9).Value = 1616  # I31: 0 -> 1616
$ws.Cells.Item(31, 10).Value = 7619.8335  # J31: 7626.9814 -> 7619.8335
$ws.Cells.Item(31, 11).Value = 1616  # K31: 0 -> 1616
$ws.Cells.Item(31, 12).Value = 7619.8335  # L31: 7626.9814 -> 7619.8335
$ws.Cells.Item(31, 13).Value = -1321  # M31: None -> -1321
$ws.Cells.Item(31, 14).Value = -8209.833500000001  # N31: -8216.981400000001 -> -8209.833500000001
# Row 34
$ws.Cells.Item(34, 8).Value = 5618.5557  # H34: 7626.9814 -> 5618.5557
$ws.Cells.Item(34, 9).Value = 1616  # I34: 0 -> 1616
$ws.Cells.Item(34, 10).Value = 7619.8335  # J34: 7626.9814 -> 7619.8335
$ws.Cells.Item(34, 11).Value = 1616  # K34: 0 -> 1616
$ws.Cells.Item(34, 12).Value = 7619.8335  # L34: 7626.9814 -> 7619.8335
$ws.Cells.Item(34, 13).Value = -1414  # M34: None -> -1414
$ws.Cells.Item(34, 14).Value = -8023.8335  # N34: -8030.9814 -> -8023.8335
# Row 107
$ws.Cells.Item(107, 8).Value = 2404927.2  # H107: 2718080.8 -> 2404927.2
$ws.Cells.Item(107, 9).Value = 3677258.5  # I107: 3907026 -> 3677258.5
$ws.Cells.Item(107, 10).Value = 1634.7778  # J107: 491.57144 -> 1634.7778
$ws.Cells.Item(107, 11).Value = 3677258.5  # K107: 3907026 -> 3677258.5
$ws.Cells.Item(107, 12).Value = 1634.7778  # L107: 491.57144 -> 1634.7778
$ws.Cells.Item(107, 13).Value = -3675338.5  # M107: -3905106 -> -3675338.5
$ws.Cells.Item(107, 14).Value = -5474.7778  # N107: -4331.57144 -> -5474.7778
# Row 140
$ws.Cells.Item(140, 8).Value = 63238  # H140: 64262.5 -> 63238
$ws.Cells.Item(140, 10).Value = 63238  # J140: 64262.5 -> 63238
$ws.Cells.Item(140, 12).Value = 63238  # L140: 64262.5 -> 63238
$ws.Cells.Item(140, 14).Value = -73598  # N140: -74622.5 -> -73598

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Cells.Item(22, 8).Value = 1020.08  # H22: 934.25 -> 1020.08
$ws.Cells.Item(22, 9).Value = 900  # I22: 650 -> 900
$ws.Cells.Item(22, 10).Value = 1050.1  # J22: 999.8461 -> 1050.1
$ws.Cells.Item(22, 11).Value = 2700  # K22: 1950 -> 2700
$ws.Cells.Item(22, 12).Value = 3150.3  # L22: 2999.5383 -> 3150.3
$ws.Cells.Item(22, 13).Value = -2531  # M22: -1781 -> -2531
$ws.Cells.Item(22, 14).Value = -3488.3  # N22: -3337.5383 -> -3488.3
# Row 27
$ws.Cells.Item(27, 8).Value = 1020.08  # H27: 934.25 -> 1020.08
$ws.Cells.Item(27, 9).Value = 900  # I27: 650 -> 900
$ws.Cells.Item(27, 10).Value = 1050.1  # J27: 999.8461 -> 1050.1
$ws.Cells.Item(27, 11).Value = 2700  # K27: 1950 -> 2700
$ws.Cells.Item(27, 12).Value = 3150.3  # L27: 2999.5383 -> 3150.3
$ws.Cells.Item(27, 13).Value = -2598  # M27: -1848 -> -2598
$ws.Cells.Item(27, 14).Value = -3354.3  # N27: -3203.5383 -> -3354.3
# Row 49
$ws.Cells.Item(49, 8).Value = 7728.8125  # H49: 7709 -> 7728.8125
$ws.Cells.Item(49, 10).Value = 7728.8125  # J49: 7709 -> 7728.8125
$ws.Cells.Item(49, 12).Value = 23186.4375  # L49: 23127 -> 23186.4375
$ws.Cells.Item(49, 14).Value = -23498.4375  # N49: -23439 -> -23498.4375
# Row 103
$ws.Cells.Item(103, 8).Value = 4584.1177  # H103: 2304 -> 4584.1177
$ws.Cells.Item(103, 9).Value = 510  # I103: 608.3333 -> 510
$ws.Cells.Item(103, 10).Value = 6281.6665  # J103: 3999.6667 -> 6281.6665
$ws.Cells.Item(103, 11).Value = 1530  # K103: 1824.9999 -> 1530
$ws.Cells.Item(103, 12).Value = 18844.9995  # L103: 11999.0001 -> 18844.9995
$ws.Cells.Item(103, 13).Value = -651  # M103: -945.9999 -> -651
$ws.Cells.Item(103, 14).Value = -20602.9995  # N103: -13757.0001 -> -20602.9995
# Row 113
$ws.Cells.Item(113, 8).Value = 572.15216  # H113: 573.75555 -> 572.15216
$ws.Cells.Item(113, 10).Value = 679.4545000000001  # J113: 697.4 -> 679.4545000000001
$ws.Cells.Item(113, 12).Value = 2038.3635  # L113: 2092.2 -> 2038.3635
$ws.Cells.Item(113, 14).Value = -6378.3635  # N113: -6432.2 -> -6378.3635
# Row 118
$ws.Cells.Item(118, 8).Value = 3512.3555  # H118: 4261.75 -> 3512.3555
$ws.Cells.Item(118, 9).Value = 1599.6666  # I118: 2301.2856 -> 1599.6666
$ws.Cells.Item(118, 10).Value = 3806.6155  # J118: 7006.4 -> 3806.6155
$ws.Cells.Item(118, 11).Value = 4798.9998  # K118: 6903.8568 -> 4798.9998
$ws.Cells.Item(118, 12).Value = 11419.8465  # L118: 21019.2 -> 11419.8465
$ws.Cells.Item(118, 13).Value = -3555.9998  # M118: -5660.8568 -> -3555.9998
$ws.Cells.Item(118, 14).Value = -13905.8465  # N118: -23505.2 -> -13905.8465
# Row 131
$ws.Cells.Item(131, 8).Value = 3921.6  # H131: 3920.75 -> 3921.6
$ws.Cells.Item(131, 10).Value = 4403.2573  # J131: 4402.2856 -> 4403.2573
$ws.Cells.Item(131, 12).Value = 13209.7719  # L131: 13206.8568 -> 13209.7719
$ws.Cells.Item(131, 14).Value = -23289.7719  # N131: -23286.8568 -> -23289.7719
# Row 137
$ws.Cells.Item(137, 8).Value = 31811.5  # H137: 36146.605 -> 31811.5
$ws.Cells.Item(137, 10).Value = 62629.35  # J137: 87391.664 -> 62629.35
$ws.Cells.Item(137, 12).Value = 187888.05  # L137: 262174.992 -> 187888.05
$ws.Cells.Item(137, 14).Value = -198088.05  # N137: -272374.992 -> -198088.05

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 1545  # H122: 1605.4 -> 1545
$ws.Cells.Item(122, 9).Value = 990  # I122: 1459.8 -> 990
$ws.Cells.Item(122, 10).Value = 2100  # J122: 1751 -> 2100
$ws.Cells.Item(122, 11).Value = 2970  # K122: 4379.4 -> 2970
$ws.Cells.Item(122, 12).Value = 6300  # L122: 5253 -> 6300
$ws.Cells.Item(122, 13).Value = -520  # M122: -1929.4 -> -520
$ws.Cells.Item(122, 14).Value = -11200  # N122: -10153 -> -11200
# Row 132
$ws.Cells.Item(132, 8).Value = 32263384  # H132: 41673420 -> 32263384
$ws.Cells.Item(132, 9).Value = 45460540  # I132: 66675270 -> 45460540
$ws.Cells.Item(132, 11).Value = 136381620  # K132: 200025810 -> 136381620
$ws.Cells.Item(132, 13).Value = -136379090  # M132: -200023280 -> -136379090
# Row 133
$ws.Cells.Item(133, 8).Value = 65780  # H133: 60778 -> 65780
$ws.Cells.Item(133, 10).Value = 65780  # J133: 60778 -> 65780
$ws.Cells.Item(133, 12).Value = 65780  # L133: 60778 -> 65780
$ws.Cells.Item(133, 14).Value = -75900  # N133: -70898 -> -75900
# Row 141
$ws.Cells.Item(141, 8).Value = 80000  # H141: 21374.625 -> 80000
$ws.Cells.Item(141, 9).Value = 0  # I141: 20000 -> 0
$ws.Cells.Item(141, 10).Value = 80000  # J141: 22199.4 -> 80000
$ws.Cells.Item(141, 11).Value = 0  # K141: 20000 -> 0
$ws.Cells.Item(141, 12).Value = 80000  # L141: 22199.4 -> 80000
$ws.Cells.Item(141, 13).ClearContents()  # M141: -14820 -> (removed)
$ws.Cells.Item(141, 14).Value = -90360  # N141: -32559.4 -> -90360

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 6772.1177  # H22: 6423.6665 -> 6772.1177
$ws.Cells.Item(22, 10).Value = 10997.6  # J22: 10043.272 -> 10997.6
$ws.Cells.Item(22, 12).Value = 10997.6  # L22: 10043.272 -> 10997.6
$ws.Cells.Item(22, 14).Value = -11587.6  # N22: -10633.272 -> -11587.6
# Row 27
$ws.Cells.Item(27, 8).Value = 6772.1177  # H27: 6423.6665 -> 6772.1177
$ws.Cells.Item(27, 10).Value = 10997.6  # J27: 10043.272 -> 10997.6
$ws.Cells.Item(27, 12).Value = 10997.6  # L27: 10043.272 -> 10997.6
$ws.Cells.Item(27, 14).Value = -11211.6  # N27: -10257.272 -> -11211.6
# Row 63
$ws.Cells.Item(63, 8).Value = 49395.715  # H63: 0 -> 49395.715
$ws.Cells.Item(63, 10).Value = 49395.715  # J63: 0 -> 49395.715
$ws.Cells.Item(63, 12).Value = 49395.715  # L63: 0 -> 49395.715
$ws.Cells.Item(63, 14).Value = -50893.715  # N63: None -> -50893.715
# Row 66
$ws.Cells.Item(66, 8).Value = 49395.715  # H66: 0 -> 49395.715
$ws.Cells.Item(66, 10).Value = 49395.715  # J66: 0 -> 49395.715
$ws.Cells.Item(66, 12).Value = 148187.145  # L66: 0 -> 148187.145
$ws.Cells.Item(66, 14).Value = -155675.145  # N66: None -> -155675.145
# Row 82
$ws.Cells.Item(82, 8).Value = 125001990  # H82: 1481.5 -> 125001990
$ws.Cells.Item(82, 9).Value = 250001000  # I82: 1534 -> 250001000
$ws.Cells.Item(82, 10).Value = 2990  # J82: 1450 -> 2990
$ws.Cells.Item(82, 11).Value = 250001000  # K82: 1534 -> 250001000
$ws.Cells.Item(82, 12).Value = 2990  # L82: 1450 -> 2990
$ws.Cells.Item(82, 13).Value = -250000639  # M82: -1173 -> -250000639
$ws.Cells.Item(82, 14).Value = -3712  # N82: -2172 -> -3712
# Row 85
$ws.Cells.Item(85, 8).Value = 125001990  # H85: 1481.5 -> 125001990
$ws.Cells.Item(85, 9).Value = 250001000  # I85: 1534 -> 250001000
$ws.Cells.Item(85, 10).Value = 2990  # J85: 1450 -> 2990
$ws.Cells.Item(85, 11).Value = 250001000  # K85: 1534 -> 250001000
$ws.Cells.Item(85, 12).Value = 2990  # L85: 1450 -> 2990
$ws.Cells.Item(85, 13).Value = -249999752  # M85: -286 -> -249999752
$ws.Cells.Item(85, 14).Value = -5486  # N85: -3946 -> -5486
# Row 136
$ws.Cells.Item(136, 8).Value = 5749293  # H136: 6412698 -> 5749293
$ws.Cells.Item(136, 9).Value = 1974.3  # I136: 2069.8235 -> 1974.3
$ws.Cells.Item(136, 10).Value = 18521112  # J136: 18521662 -> 18521112
$ws.Cells.Item(136, 11).Value = 5922.9  # K136: 6209.470499999999 -> 5922.9
$ws.Cells.Item(136, 12).Value = 55563336  # L136: 55564986 -> 55563336
$ws.Cells.Item(136, 13).Value = -3372.9  # M136: -3659.470499999999 -> -3372.9
$ws.Cells.Item(136, 14).Value = -55568436  # N136: -55570086 -> -55568436
# Row 140
$ws.Cells.Item(140, 8).Value = 65246.637  # H140: 68610.625 -> 65246.637
$ws.Cells.Item(140, 10).Value = 65246.637  # J140: 68610.625 -> 65246.637
$ws.Cells.Item(140, 12).Value = 65246.637  # L140: 68610.625 -> 65246.637
$ws.Cells.Item(140, 14).Value = -75606.637  # N140: -78970.625 -> -75606.637
# Row 141
$ws.Cells.Item(141, 8).Value = 61473.2  # H141: 60690.055 -> 61473.2
$ws.Cells.Item(141, 10).Value = 61473.2  # J141: 60690.055 -> 61473.2
$ws.Cells.Item(141, 12).Value = 61473.2  # L141: 60690.055 -> 61473.2
$ws.Cells.Item(141, 14).Value = -71833.2  # N141: -71050.05499999999 -> -71833.2

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 4355459  # H132: 4421463 -> 4355459
$ws.Cells.Item(132, 9).Value = 2305.2856  # I132: 2240.8 -> 2305.2856
$ws.Cells.Item(132, 10).Value = 11668757  # J132: 13891224 -> 11668757
$ws.Cells.Item(132, 11).Value = 6915.8568  # K132: 6722.400000000001 -> 6915.8568
$ws.Cells.Item(132, 12).Value = 35006271  # L132: 41673672 -> 35006271
$ws.Cells.Item(132, 13).Value = -4385.8568  # M132: -4192.400000000001 -> -4385.8568
$ws.Cells.Item(132, 14).Value = -35011331  # N132: -41678732 -> -35011331
# Row 141
$ws.Cells.Item(141, 8).Value = 69522.5  # H141: 69723.81 -> 69522.5
$ws.Cells.Item(141, 10).Value = 69522.5  # J141: 69723.81 -> 69522.5
$ws.Cells.Item(141, 12).Value = 69522.5  # L141: 69723.81 -> 69522.5
$ws.Cells.Item(141, 14).Value = -79882.5  # N141: -80083.81 -> -79882.5
